$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '67.217.32'
$cell.ClearFormats()

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '  -3.58%  '
$cell.ClearFormats()

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.487.35'
$cell.ClearFormats()

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '  -5.24%  '
$cell.ClearFormats()

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '  -0.07%  '
$cell.ClearFormats()

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '608.08'
$cell.ClearFormats()

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '  -6.65%  '
$cell.ClearFormats()

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '148.23'
$cell.ClearFormats()

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '  -7.93%  '
$cell.ClearFormats()

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '3.485.91'
$cell.ClearFormats()

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '  -5.18%  '
$cell.ClearFormats()

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '  -0.02%  '
$cell.ClearFormats()

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '  -3.84%  '
$cell.ClearFormats()

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.138'
$cell.ClearFormats()

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '  -5.10%  '
$cell.ClearFormats()

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '6.92'
$cell.ClearFormats()

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '  -3.26%  '
$cell.ClearFormats()

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.422'
$cell.ClearFormats()

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '  -4.49%  '
$cell.ClearFormats()

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '  -6.56%  '
$cell.ClearFormats()

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '4.074.89'
$cell.ClearFormats()

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '  -5.29%  '
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '31.38'
$cell.ClearFormats()

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '  -3.89%  '
$cell.ClearFormats()

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '3.484.91'
$cell.ClearFormats()

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '  -4.89%  '
$cell.ClearFormats()

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '67.145.67'
$cell.ClearFormats()

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '  -3.73%  '
$cell.ClearFormats()

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '  -0.94%  '
$cell.ClearFormats()

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '  -1.65%  '
$cell.ClearFormats()

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '14.99'
$cell.ClearFormats()

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '  -5.84%  '
$cell.ClearFormats()

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '446.44'
$cell.ClearFormats()

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '  -5.15%  '
$cell.ClearFormats()

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '9.03'
$cell.ClearFormats()

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '  -12.83%  '
$cell.ClearFormats()

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.622'
$cell.ClearFormats()

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '  -4.89%  '
$cell.ClearFormats()

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '77.01'
$cell.ClearFormats()

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '  -3.43%  '
$cell.ClearFormats()

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '  +0.12%  '
$cell.ClearFormats()

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '  +0.01%  '
$cell.ClearFormats()

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '3.625.71'
$cell.ClearFormats()

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '  -5.26%  '
$cell.ClearFormats()

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '10.11'
$cell.ClearFormats()

$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '  -8.94%  '
$cell.ClearFormats()

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '8.28'
$cell.ClearFormats()

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = '  -5.30%  '
$cell.ClearFormats()

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '2.53'
$cell.ClearFormats()

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = '  -4.59%  '
$cell.ClearFormats()

$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = '  -7.52%  '
$cell.ClearFormats()

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = '  +0.24%  '
$cell.ClearFormats()

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.165'
$cell.ClearFormats()

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = '  -0.79%  '
$cell.ClearFormats()

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '25.73'
$cell.ClearFormats()

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = '  -3.74%  '
$cell.ClearFormats()

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = '  -4.95%  '
$cell.ClearFormats()

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '1.85'
$cell.ClearFormats()

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = '  -7.15%  '
$cell.ClearFormats()

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '3.480.51'
$cell.ClearFormats()

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = '  -5.30%  '
$cell.ClearFormats()

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '7.97'
$cell.ClearFormats()

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '  -4.75%  '
$cell.ClearFormats()

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '  +0.06%  '
$cell.ClearFormats()

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '  -0.12%  '
$cell.ClearFormats()

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '  +0.60%  '
$cell.ClearFormats()

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '170.66'
$cell.ClearFormats()

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '  -4.25%  '
$cell.ClearFormats()

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.0871'
$cell.ClearFormats()

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '  -2.57%  '
$cell.ClearFormats()

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '5.42'
$cell.ClearFormats()

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '  -7.64%  '
$cell.ClearFormats()

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.882'
$cell.ClearFormats()

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '  -5.12%  '
$cell.ClearFormats()

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '  -3.08%  '
$cell.ClearFormats()

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '  +1.63%  '
$cell.ClearFormats()

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '26.49'
$cell.ClearFormats()

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '  -8.85%  '
$cell.ClearFormats()

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '  -8.93%  '
$cell.ClearFormats()

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '7.52'
$cell.ClearFormats()

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '  -4.27%  '
$cell.ClearFormats()

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.01'
$cell.ClearFormats()

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '  -3.91%  '
$cell.ClearFormats()
